$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("G9","G10","G11","G13","G14","G19","H19","G21","H21")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C8").Value = 66
$ws.Range("C9").Value = 30
$ws.Range("G9").Value = "7680.00"
$ws.Range("C10").Value = 20
$ws.Range("G10").Value = "9440.00"
$ws.Range("C11").Value = 88
$ws.Range("G11").Value = "58256.00"
$ws.Range("C12").Value = 98
$ws.Range("C13").Value = 20
$ws.Range("G13").Value = "2720.00"
$ws.Range("C14").Value = 94
$ws.Range("G14").Value = "2162.00"
$ws.Range("C15").Value = 10
$ws.Range("C16").Value = 100
$ws.Range("C17").Value = 81
$ws.Range("G19").Value = "80258.00"
$ws.Range("H19").Value = "80258.00"
$ws.Range("G21").Value = "80258.00"
$ws.Range("H21").Value = "80258.00"
